# "product control system completed"
# Turns the plain Stock sheet into a small product-control table:
#   - inserts a new "PRO_N0" id column in front of TÜR/ADET
#   - adds two more products (göt, ruvi) with their ADET counts
#   - makes the Stock sheet the active tab/sheet (was Sayfa1)

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Stock")

# --- Stock sheet: insert a new first column for the product number ---
$ws2.Range("A1").EntireColumn.Insert()

# Header row
$ws2.Cells.Item(1,1).Value = "PRO_N0"
# B1 "TÜR" and C1 "ADET" already shifted into place by the column insert.

# Existing rows get sequential product numbers in the new column A
$ws2.Cells.Item(2,1).Value = 1
$ws2.Cells.Item(3,1).Value = 2

# Two new products appended as rows 4 and 5
$ws2.Cells.Item(4,1).Value = 3
$ws2.Cells.Item(4,2).Value = "göt"

$ws2.Cells.Item(5,1).Value = 4
$ws2.Cells.Item(5,2).Value = "ruvi"

# ADET (quantity) column - keep these as text, matching the source data
$qty = $ws2.Range("C2:C5")
$qty.NumberFormat = "@"
$ws2.Cells.Item(2,3).Value = "1"
$ws2.Cells.Item(3,3).Value = "1"
$ws2.Cells.Item(4,3).Value = "12"
$ws2.Cells.Item(5,3).Value = "2"

# --- Make Stock the active sheet/tab, with its working range selected ---
$ws2.Activate()
$ws2.Range("A1:I10").Select()

Write-Output "Stock sheet populated and activated"
